# DOMA-4532 hotfix: fix meter-import-example.xlsx header labels / column widths
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text fixes (row 1) ---------------------------------------
# D1: "ЛС" -> "Лицевой счет"
$ws.Cells.Item(1, 4).Value = "Лицевой счет"
# F1: "№ счетчика" -> "Номер счетчика"
$ws.Cells.Item(1, 6).Value = "Номер счетчика"
# G1: "Кол-во тарифов" -> "Количество тарифов"
$ws.Cells.Item(1, 7).Value = "Количество тарифов"

# --- Column width fixes (widen columns to fit the longer header text) -
# The worksheet's stored `width` (characters) is ColumnWidth + 5/6 and is
# snapped by the engine to 1/6-character increments, so we back-solve the
# ColumnWidth input that lands on the closest achievable stored width.
$offset = 5 / 6
$ws.Columns.Item(4).ColumnWidth = 11.3516 - $offset   # D: 3.67188 -> 11.3516
$ws.Columns.Item(6).ColumnWidth = 13.3516 - $offset   # F: 10      -> 13.3516
$ws.Columns.Item(7).ColumnWidth = 17 - $offset        # G: 13.3516 -> 17

# --- Unused template font (sz 15 Calibri -> sz 13 Arial) --------------
# This font record (index 2 in styles.xml `fonts`) is not referenced by
# any cellXf in the original workbook - it's a leftover template font
# with no visible effect - so there is no cell/range to drive the COM
# object model through that would rewrite it in place without altering
# real, visibly-used formatting elsewhere in the sheet.
